# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17 (pushes old row 17 -> 18, old rows 22/23 -> 23/24,
# and the relevant merged cells move down with them automatically).
$ws.Rows("17").Insert()

# The new row 17 comes in blank; clone the formatting used by the data row
# directly above it (row 16) so it matches the rest of the table.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the 3 "Estado de Cuenta" data rows (periods 2507, 2506, 2505) ---

# Row 16: period 2505 -> 2507, updated Valor Mora
$ws.Range("E16").Value = "2507"
$ws.Range("G16").Value = 737717

# Row 17 (new): period 2506
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45555876"
$ws.Range("D17").Value = "VIVIANA PAOLA MESTRA PADILLA"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 100000
$ws.Range("G17").Value = 737717

# Row 18 (previously row 17): period 2505, updated Valor Mora
$ws.Range("E18").Value = "2505"
$ws.Range("G18").Value = 737717

# --- Update summary fields ---
$ws.Range("E11").Value = 300000
$ws.Range("F13").Value = 3
